# Regenerate the "K" (strikeouts) column (column G) of the save_data sheet
# with corrected values, replacing the previous "Strike#" based figures.
# Commit: "regen save_data to use K instead of Strike#, regen std/mean,
#          calc and write s_vals"

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Mapping of row number -> new value for column G ("K")
$kValues = @{
    2  = 0
    3  = 1
    4  = 0
    5  = 1
    6  = 1
    7  = 1
    8  = 1
    9  = 2
    10 = 0
    11 = 2
    12 = 1
    13 = 0
    14 = 1
    15 = 0
    16 = 1
    17 = 1
    18 = 2
    19 = 1
    20 = 0
    21 = 0
    22 = 1
    23 = 2
    24 = 0
    25 = 1
    26 = 0
    27 = 0
    28 = 1
    29 = 0
    30 = 0
    31 = 0
    32 = 1
    33 = 1
    34 = 1
    35 = 2
    36 = 3
    37 = 0
    38 = 2
    39 = 1
    40 = 0
    41 = 1
    42 = 0
    43 = 1
    44 = 1
    45 = 2
    46 = 1
    47 = 1
    48 = 3
    49 = 2
    50 = 0
    51 = 1
    52 = 2
    53 = 0
    54 = 1
    55 = 2
    56 = 1
    57 = 2
    58 = 0
    59 = 0
    60 = 0
}

foreach ($row in $kValues.Keys) {
    $ws.Cells.Item($row, 7).Value = $kValues[$row]
}
